# Presentation04 - "add initial presentation 4 draft"
#
# Fills in the titles for the three "Interesting Things" slides, the
# "Biggest Challenge" title, and the two still-empty first body
# placeholders (one on the "If You Could Do It Again..." slide, one on
# the "R Shiny" slide) with the authors' draft content.

$p = $ppt.ActivePresentation

# Slide 2: "Interesting Things: 1" -> "Interesting Things: Nonparametric fails"
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Interesting Things: Nonparametric fails"

# Slide 3: "Interesting Things: 2" -> three runs, splitting out "xgboost"
# (flagged by the spell checker in the original file) as its own run.
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Text = "Interesting Things: xgboost for the win"
$tr3.Characters(21, 7).Text = "xgboost"

# Slide 4: "Interesting Things: 3" -> "Interesting Things: No one size fits all"
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Interesting Things: No one size fits all"

# Slide 5: "Biggest Challenge " -> "Biggest Challenge: API and fees "
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Biggest Challenge: API and fees "

# Slide 6: first body placeholder was empty -> "Precipitation data"
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Precipitation data"

# Slide 7: first body placeholder was empty -> "Jason and Aubrey's interesting things"
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(2).TextFrame.TextRange.Text = "Jason and Aubrey’s interesting things"
